$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
for ($r = 1; $r -le 19; $r++) {
    $v = $ws2.Cells.Item($r,1).Value2
    $s = $ws2.Cells.Item($r,1).Style.Name()
    Write-Host $r ":" $s ":" $v
}
